# "Report summary variables by group"
#
# The "Page one" sheet is restructured:
#   - A new "Objectives" summary table is added at the top (rows 1-4).
#   - The three existing report tables ("Birth to last vote for David"
#     (renamed/repurposed as the plain "Birth to last vote" group table),
#     "Birth to last vote", and "Last name begins with C") are each moved
#     one column to the left (column A instead of column B) and their
#     header rows change from Value/Objective/Objective met to a
#     Group/Average/Minimum/Maximum (or Group/Sum, Group/Average) layout.
#   - The per-row red/green "objective met" indicator is dropped from the
#     group tables, so the green "Yes" style is no longer used anywhere in
#     the workbook.
#
# xlPasteFormats = -4122 (captured as a constant below since this host has
# no $xlPasteFormats / Excel enum constants predefined).

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page one")

# ---- break the old merges so we can freely restyle/move content ----
$ws.Range("B1:D1").UnMerge()
$ws.Range("B6:D6").UnMerge()
$ws.Range("B11:D11").UnMerge()

# ---- copy the 4 existing direct formats onto their new homes first      ----
# ---- (reuses the workbook's existing style entries instead of creating  ----
# ---- new ones; do this before any values change so the source cells     ----
# ---- still carry their original formatting)                             ----

# style "title" (bold, grey fill, merged banner rows) - currently on B1
$ws.Range("B1").Copy()
$ws.Range("A1:B1").PasteSpecial($xlPasteFormats)
$ws.Range("A6:D6").PasteSpecial($xlPasteFormats)
$ws.Range("A11:B11").PasteSpecial($xlPasteFormats)
$ws.Range("A16:B16").PasteSpecial($xlPasteFormats)

# style "sub-header" (bold, light fill) - currently on B2
$ws.Range("B2").Copy()
$ws.Range("A2:B2").PasteSpecial($xlPasteFormats)
$ws.Range("A7:D7").PasteSpecial($xlPasteFormats)
$ws.Range("A12:B12").PasteSpecial($xlPasteFormats)
$ws.Range("A17:B17").PasteSpecial($xlPasteFormats)

# style "data" (plain, bordered) - currently on B3
$ws.Range("B3").Copy()
$ws.Range("A3:A4").PasteSpecial($xlPasteFormats)
$ws.Range("A8:D8").PasteSpecial($xlPasteFormats)
$ws.Range("A13:B13").PasteSpecial($xlPasteFormats)
$ws.Range("A18:B18").PasteSpecial($xlPasteFormats)

# style "red / not met" (bold red) - currently on D3
$ws.Range("D3").Copy()
$ws.Range("B3:B4").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---- clear every cell that is no longer part of the new layout ----
$ws.Range("C1:D1").Clear()
$ws.Range("C2:D2").Clear()
$ws.Range("C3:D3").Clear()
$ws.Range("C11:D11").Clear()
$ws.Range("C12:D12").Clear()
$ws.Range("C13:D13").Clear()
$ws.Range("A14:D15").Clear()

# ---- write the new cell values ----

# Objectives (rows 1-4, new table)
$ws.Range("A1").Value = "Objectives"
$ws.Range("A2").Value = "Objective"
$ws.Range("B2").Value = "Met"
$ws.Range("A3").Value = "Birth to last vote average less than 100"
$ws.Range("B3").Value = "No"
$ws.Range("A4").Value = "More 2 or more people with last name beginning with C"
$ws.Range("B4").Value = "No"

# Birth to last vote (rows 6-8)
$ws.Range("A6").Value = "Birth to last vote"
$ws.Range("A7").Value = "Group"
$ws.Range("B7").Value = "Average"
$ws.Range("C7").Value = "Minimum"
$ws.Range("D7").Value = "Maximum"
$ws.Range("A8").Value = "All observations"
$ws.Range("B8").Value = 13810.3333333333
$ws.Range("C8").Value = 9156
$ws.Range("D8").Value = 21127

# Last name begins with C (rows 11-13)
$ws.Range("A11").Value = "Last name begins with C"
$ws.Range("A12").Value = "Group"
$ws.Range("B12").Value = "Sum"
$ws.Range("A13").Value = "All observations"
$ws.Range("B13").Value = 1

# Birth to last vote for David (rows 16-18)
$ws.Range("A16").Value = "Birth to last vote for David"
$ws.Range("A17").Value = "Group"
$ws.Range("B17").Value = "Average"
$ws.Range("A18").Value = "First name David"
$ws.Range("B18").Value = 11148

# ---- re-merge the title banner cells ----
$ws.Range("A1:B1").Merge()
$ws.Range("A6:D6").Merge()
$ws.Range("A11:B11").Merge()
$ws.Range("A16:B16").Merge()

Write-Host "Page one rebuilt"
